# Insert a new weekly price record at row 185 ("Hortaliza" / Papa,
# Terminal Hortofrutícola Agro Chillán). Inserting the row shifts every
# existing record from row 185 onward down by one (old row 185 -> 186,
# ..., old row 223 -> 224), growing the used range from A1:R223 to
# A1:R224 — matching the diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 185..223 down to 186..224, duplicating formatting from the
# row above (Excel's default Insert behaviour) so D185 keeps its date
# style.
$ws.Rows.Item(185).Insert()

# Fill in the newly inserted row with the new record's data.
$ws.Cells.Item(185, 1).Value  = 7
$ws.Cells.Item(185, 2).Value  = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(185, 3).Value  = "Ñuble"
$ws.Cells.Item(185, 4).Value  = 44511
$ws.Cells.Item(185, 5).Value  = 16
$ws.Cells.Item(185, 6).Value  = 100114001
$ws.Cells.Item(185, 7).Value  = "Papa"
$ws.Cells.Item(185, 8).Value  = "Patagonia"
$ws.Cells.Item(185, 9).Value  = "1a (guarda)"
$ws.Cells.Item(185, 10).Value = 160
$ws.Cells.Item(185, 11).Value = 7000
$ws.Cells.Item(185, 12).Value = 8000
$ws.Cells.Item(185, 13).Value = 7500
$ws.Cells.Item(185, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(185, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(185, 16).Value = 300
$ws.Cells.Item(185, 17).Value = 25
$ws.Cells.Item(185, 18).Value = "Hortaliza"
